$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM numbers were regenerated. The two sending
# clusters that used to populate rows 2-7 ("ECs" and "MuSCs") collapse
# to a single refreshed "MuSCs" block: the old "ECs -> *" rows (2-4) are
# replaced by recalculated "MuSCs -> *" values, and the old, now
# redundant "MuSCs -> *" rows (5-7) are removed outright, so the sheet's
# used range shrinks from A1:T7 down to A1:T4.
$ws.Rows("5:7").Delete()

# Row 2: MuSCs -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.0005903333333333333
$ws.Range("H2").Value = 0.001771
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 0.002215484005777778
$ws.Range("R2").Value = 0.019939356052
$ws.Range("S2").Value = 0.6855621274031838
$ws.Range("T2").Value = 0.6855621274031838

# Row 3: MuSCs -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.0005903333333333333
$ws.Range("H3").Value = 0.001771
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.296447666666667
$ws.Range("N3").Value = 3.889343
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 0.0007653362725555555
$ws.Range("R3").Value = 0.006888026453000001
$ws.Range("S3").Value = 0.2368266084628361
$ws.Range("T3").Value = 0.2368266084628362

# Row 4: MuSCs -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.0005903333333333333
$ws.Range("H4").Value = 0.001771
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 0.0002508109877777778
$ws.Range("R4").Value = 0.00225729889
$ws.Range("S4").Value = 0.07761126413398003
$ws.Range("T4").Value = 0.07761126413398005
